$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.892.61"
$ws.Range("E2").Value = "  +4.21%  "

$ws.Range("D3").Value = "'3.249.40"
$ws.Range("E3").Value = "  +3.12%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "'545.15"
$ws.Range("E5").Value = "  +3.29%  "

$ws.Range("D6").Value = "'147.20"
$ws.Range("E6").Value = "  +5.49%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.528"
$ws.Range("E8").Value = "  -0.32%  "

$ws.Range("D9").Value = "'7.37"
$ws.Range("E9").Value = "  +1.51%  "

$ws.Range("E10").Value = "  +3.32%  "

$ws.Range("D11").Value = "'0.432"
$ws.Range("E11").Value = "  -1.30%  "

$ws.Range("D12").Value = "'3.808.58"
$ws.Range("E12").Value = "  +3.01%  "

$ws.Range("E13").Value = "  -1.84%  "

$ws.Range("D14").Value = "'26.26"
$ws.Range("E14").Value = "  +2.66%  "

$ws.Range("E15").Value = "  +3.43%  "

$ws.Range("D16").Value = "'60.814.41"
$ws.Range("E16").Value = "  +3.82%  "

$ws.Range("D17").Value = "'3.261.13"
$ws.Range("E17").Value = "  +3.04%  "

$ws.Range("D18").Value = "'6.31"
$ws.Range("E18").Value = "  +1.87%  "

$ws.Range("D19").Value = "'13.42"
$ws.Range("E19").Value = "  +4.09%  "

$ws.Range("D20").Value = "'8.39"
$ws.Range("E20").Value = "  +3.67%  "

$ws.Range("D21").Value = "'377.88"
$ws.Range("E21").Value = "  +1.39%  "

$ws.Range("E22").Value = "  +0.23%  "

$ws.Range("D23").Value = "'0.530"
$ws.Range("E23").Value = "  +0.77%  "

$ws.Range("D24").Value = "'69.93"
$ws.Range("E24").Value = "  +0.52%  "

$ws.Range("D25").Value = "'0.170"
$ws.Range("E25").Value = "  +1.68%  "

$ws.Range("D26").Value = "'8.64"
$ws.Range("E26").Value = "  +2.88%  "

$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("D28").Value = "'0.0₃0911"
$ws.Range("E28").Value = "  +7.45%  "

$ws.Range("E29").Value = "  +3.08%  "

$ws.Range("D30").Value = "'22.58"
$ws.Range("E30").Value = "  +1.11%  "

$ws.Range("D31").Value = "'6.18"
$ws.Range("E31").Value = "  +3.52%  "

$ws.Range("D32").Value = "'5.41"
$ws.Range("E32").Value = "  +5.74%  "

$ws.Range("E33").Value = "  +7.70%  "

$ws.Range("D34").Value = "'6.64"
$ws.Range("E34").Value = "  +5.69%  "

$ws.Range("D35").Value = "'159.20"
$ws.Range("E35").Value = "  +1.45%  "

$ws.Range("E36").Value = "  +7.94%  "

$ws.Range("D37").Value = "'26.37"
$ws.Range("E37").Value = "  +5.97%  "

$ws.Range("D38").Value = "'2.817.58"
$ws.Range("E38").Value = "  +4.67%  "

$ws.Range("E39").Value = "  +9.24%  "

$ws.Range("E40").Value = "  +4.39%  "

$ws.Range("D41").Value = "'1.72"
$ws.Range("E41").Value = "  +2.71%  "

$ws.Range("D42").Value = "'4.29"
$ws.Range("E42").Value = "  +1.01%  "

$ws.Range("D43").Value = "'40.01"
$ws.Range("E43").Value = "  +2.62%  "

$ws.Range("D44").Value = "'0.730"
$ws.Range("E44").Value = "  +1.25%  "

$ws.Range("D45").Value = "'3.289.06"
$ws.Range("E45").Value = "  +2.89%  "

$ws.Range("E46").Value = "  +3.16%  "

$ws.Range("D47").Value = "'1.00"
$ws.Range("E47").Value = "  +2.86%  "

$ws.Range("D48").Value = "'21.37"
$ws.Range("E48").Value = "  +7.09%  "

$ws.Range("D49").Value = "'6.25"
$ws.Range("E49").Value = "  +1.10%  "

$ws.Range("D50").Value = "'0.801"
$ws.Range("E50").Value = "  +7.70%  "

$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'274.44"
$ws.Range("E51").Value = "  +7.22%  "
